# TC23_Verify_UserRegistration.xlsx - "Update in User registration"
#
# For every ENTERTEXT / ENTER_RANDOM_VALUE step (FName, LName, Email, Pass,
# ConfirmPass, Zip) a new "CLICK_PRE_ENTERTEXT" step is inserted directly
# above it, clicking on the same target Object before the text gets typed
# into it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# (row-to-insert-before, Object name used by the pair of rows)
# Processed bottom-to-top so earlier row numbers remain valid while we insert.
$targets = @(
    @{ Row = 11; Object = "RegistrationZip" },
    @{ Row = 10; Object = "RegistrationConfirmPass" },
    @{ Row = 9;  Object = "RegistrationPass" },
    @{ Row = 8;  Object = "RegistrationEmail" },
    @{ Row = 7;  Object = "RegistrationLname" },
    @{ Row = 6;  Object = "RegistrationFName" }
)

foreach ($t in $targets) {
    $row = $t.Row

    # Insert a blank row above the existing ENTERTEXT row, shifting it (and
    # everything below) down by one.
    $ws.Rows.Item($row).Insert()

    # A column stays blank.
    $ws.Cells.Item($row, 2).Value2 = "CLICK_PRE_ENTERTEXT"
    $ws.Cells.Item($row, 3).Value2 = $t.Object
    $ws.Cells.Item($row, 4).Value2 = "CSS"
    # E column is left empty for the new CLICK_PRE_ENTERTEXT row.

    # Match the bordered-cell look used throughout the table.
    $newRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 5))
    $newRange.Borders.LineStyle = 1
}

# Sheet now spans through row 29 instead of row 23.
$ws.Range("A1:E29").Select()

# Reflect the view state recorded in the saved workbook: scrolled down a bit,
# with the last (blank SCROLL helper) row selected.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("A28:XFD28").Select()
$ws.Range("A28").Activate()
